# Updated cryptos list with GitHub Actions
# Applies the latest Coinranking snapshot (price + 1h volume change, plus a
# few re-ranked rows) to Sheet1's data table (rows 2-51, columns A-E).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows whose Coin name / Link / Price / Volume(1h) all changed (re-ranking).
# Columns: Row, Coin, Link, Price, Volume(1h)
$rowUpdates = @(
    @(27, "RenderToken", "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr", "10.01", "  -0.76%  "),
    @(28, "Dai", "https://coinranking.com/coin/MoTuySvg7+dai-dai", "1.00", "  -0.07%  "),
    @(36, "Hedera", "https://coinranking.com/coin/jad286TjB+hedera-hbar", "0.100", "  -0.12%  "),
    @(37, "dogwifhat", "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif", "3.43", "  +2.88%  "),
    @(38, "Kaspa", "https://coinranking.com/coin/V8GxkwWow+kaspa-kas", "0.138", "  -0.15%  "),
    @(39, "Mantle", "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt", "0.999", "  +0.16%  "),
    @(40, "Filecoin", "https://coinranking.com/coin/ymQub4fuB+filecoin-fil", "5.78", "  +0.51%  "),
    @(41, "FirstDigitalUSD", "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd", "0.999", "  -0.02%  "),
    @(42, "USDe", "https://coinranking.com/coin/exbfr2U-0+usde-usde", "1.00", "  -0.02%  "),
    @(43, "OKB", "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb", "48.11", "  +2.10%  "),
    @(44, "TheGraph", "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt", "0.301", "  +1.00%  "),
    @(45, "EnergySwap", "https://coinranking.com/coin/SbWqqTui-+energyswap-ens", "28.43", "  +12.63%  "),
    @(46, "Arweave", "https://coinranking.com/coin/7XWg41D1+arweave-ar", "42.93", "  -4.01%  "),
    @(47, "Cosmos", "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom", "8.34", "  -0.45%  "),
    @(48, "ONDO", "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo", "1.36", "  +8.88%  "),
    @(49, "Monero", "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr", "147.52", "  -0.21%  "),
    @(50, "Stacks", "https://coinranking.com/coin/mMPrMcB7+stacks-stx", "1.85", "  +0.58%  "),
    @(51, "Bittensor", "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao", "391.33", "  -0.08%  ")
)

# Rows that only got a refreshed Price (kept the same Coin/Link).
# Columns: Row, Price
$priceUpdates = @(
    @(2, "67.749.81"),
    @(3, "3.803.23"),
    @(5, "599.53"),
    @(6, "167.53"),
    @(10, "6.30"),
    @(13, "35.96"),
    @(14, "4.439.57"),
    @(15, "3.744.64"),
    @(17, "67.809.39"),
    @(20, "461.94"),
    @(22, "0.701"),
    @(23, "0.0000151"),
    @(24, "83.35"),
    @(25, "12.09"),
    @(26, "2.11"),
    @(29, "3.950.86"),
    @(31, "7.42"),
    @(33, "29.52")
)

# Rows that only got a refreshed Volume(1h) (Price unchanged).
# Columns: Row, Volume(1h)
$volumeOnlyUpdates = @(
    @(4, "  +0.05%  "),
    @(7, "  -0.03%  "),
    @(8, "  +0.35%  "),
    @(9, "  +0.95%  "),
    @(11, "  +0.06%  "),
    @(12, "  -0.65%  "),
    @(16, "  -0.38%  "),
    @(18, "  +1.61%  "),
    @(19, "  +0.56%  "),
    @(21, "  -2.99%  "),
    @(30, "  -0.51%  "),
    @(32, "  +1.75%  "),
    @(34, "  -0.04%  "),
    @(35, "  -1.54%  ")
)

# Every touched row also gets its Volume(1h) cell refreshed; map row -> new value.
$volumeForRow = @{
    2  = "  +0.12%  "
    3  = "  +0.31%  "
    5  = "  +0.76%  "
    6  = "  +0.78%  "
    10 = "  -0.96%  "
    13 = "  -0.76%  "
    14 = "  +0.37%  "
    15 = "  -1.14%  "
    17 = "  +0.23%  "
    20 = "  +0.66%  "
    22 = "  +0.38%  "
    23 = "  -0.12%  "
    24 = "  -0.35%  "
    25 = "  +1.88%  "
    26 = "  -1.29%  "
    29 = "  +0.35%  "
    31 = "  +1.61%  "
    33 = "  -1.28%  "
}

function Set-TextValue($range, [string]$value) {
    # Force text storage so numeric-looking strings (e.g. "6.30", "1.00",
    # "0.0000151") keep their exact original formatting instead of being
    # coerced into a General number (which would drop trailing zeros or
    # switch to scientific notation).
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# 1) Rows with a full Coin/Link/Price/Volume refresh (re-ranked coins).
foreach ($row in $rowUpdates) {
    $r = $row[0]
    $ws.Range("B$r").Value = $row[1]
    $ws.Range("C$r").Value = $row[2]
    Set-TextValue $ws.Range("D$r") $row[3]
    $ws.Range("E$r").Value = $row[4]
}

# 2) Rows with only a Price refresh (Volume(1h) updated alongside it).
foreach ($row in $priceUpdates) {
    $r = $row[0]
    Set-TextValue $ws.Range("D$r") $row[1]
    $ws.Range("E$r").Value = $volumeForRow[$r]
}

# 3) Rows with only a Volume(1h) refresh (Price unchanged).
foreach ($row in $volumeOnlyUpdates) {
    $r = $row[0]
    $ws.Range("E$r").Value = $row[1]
}
